# Weekly fruit/vegetable price update:
# A new weekly record (Asterix, 2021-10-05) is inserted as row 204 on the
# "Papa" sheet for "Feria Lagunitas de Puerto Montt", pushing every
# subsequent record down by one row (204-234 -> 205-235).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 204, shifting rows 204:234
# down to 205:235 (formats/styles are inherited from the row above, same
# as Excel's native Insert behavior).
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(204, 1).Value  = 4
$ws.Cells.Item(204, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(204, 3).Value  = "Los Lagos"
$ws.Cells.Item(204, 4).Value  = 44474
$ws.Cells.Item(204, 5).Value  = 10
$ws.Cells.Item(204, 6).Value  = 100114001
$ws.Cells.Item(204, 7).Value  = "Papa"
$ws.Cells.Item(204, 8).Value  = "Asterix"
$ws.Cells.Item(204, 9).Value  = "1a (guarda)"
$ws.Cells.Item(204, 10).Value = 600
$ws.Cells.Item(204, 11).Value = 9000
$ws.Cells.Item(204, 12).Value = 9000
$ws.Cells.Item(204, 13).Value = 9000
$ws.Cells.Item(204, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(204, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(204, 16).Value = 360
$ws.Cells.Item(204, 17).Value = 25
$ws.Cells.Item(204, 18).Value = "Hortaliza"
